# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the file that was just handed off
# (bc991a0c-144e-4524-9c58-0401a616da59.md, row 7 on every sheet) on each
# per-language sheet, and rolls the newest of those timestamps up into the
# Overview sheet's "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# Per-language "Latest Handoff Datetime" (column H) for bc991a0c...md (row 7)
$wsZhCn.Range("H7").Value = "2016-09-06 04:53:31"
$wsDeDe.Range("H7").Value = "2016-09-06 04:53:36"

# Overview roll-up "Latest HO Xliff Generate Date" (column G) for the same
# file/row - the max of the per-language handoff timestamps above.
$wsOverview.Range("G7").Value = "2016-09-06 04:53:36"
